$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @('Andaman and Nicobar Islands', 59, 43, 0, '26-06-2020', 12, 7, 0),
  @('Andhra Pradesh', 10884, 4988, 136, '26-06-2020', 2432, 877, 35),
  @('Arunachal Pradesh', 160, 38, 1, '26-06-2020', 25, 24, 1),
  @('Assam', 6321, 4033, 9, '26-06-2020', 1417, 994, 0),
  @('Bihar', 8473, 6441, 57, '26-06-2020', 940, 912, 5),
  @('Chandigarh', 423, 329, 6, '26-06-2020', 19, 13, 0),
  @('Chhattisgarh', 2452, 1755, 12, '26-06-2020', 411, 361, 1),
  @('Dadra and Nagar Haveli and Daman and Diu', 155, 32, 0, '26-06-2020', 87, 6, 0),
  @('Delhi', 73780, 44765, 2429, '26-06-2020', 17034, 13471, 317),
  @('Goa', 995, 335, 2, '26-06-2020', 241, 206, 2),
  @('Gujarat', 29520, 21498, 1753, '26-06-2020', 2840, 2804, 115),
  @('Haryana', 12463, 7380, 198, '26-06-2020', 2240, 2252, 49),
  @('Himachal Pradesh', 839, 477, 9, '26-06-2020', 183, 64, 1),
  @('Jammu and Kashmir', 6549, 3967, 90, '26-06-2020', 715, 631, 9),
  @('Jharkhand', 2262, 1605, 12, '26-06-2020', 297, 270, 1),
  @('Karnataka', 10560, 6670, 170, '26-06-2020', 1863, 1279, 38),
  @('Kerala', 3726, 1943, 22, '26-06-2020', 687, 375, 1),
  @('Ladakh', 941, 358, 1, '26-06-2020', 105, 241, 0),
  @('Madhya Pradesh', 12596, 9619, 542, '26-06-2020', 872, 739, 41),
  @('Maharashtra', 147741, 77453, 6931, '26-06-2020', 19536, 13300, 947),
  @('Manipur', 1056, 354, 0, '26-06-2020', 279, 122, 0),
  @('Meghalaya', 46, 42, 1, '26-06-2020', 2, 9, 0),
  @('Mizoram', 145, 30, 0, '26-06-2020', 5, 21, 0),
  @('Nagaland', 355, 160, 0, '26-06-2020', 154, 22, 0),
  @('Odisha', 5962, 4291, 17, '26-06-2020', 1106, 757, 5),
  @('Puducherry', 502, 187, 9, '26-06-2020', 216, 69, 2),
  @('Punjab', 4769, 3192, 120, '26-06-2020', 817, 514, 22),
  @('Rajasthan', 16296, 12840, 379, '26-06-2020', 1760, 1566, 42),
  @('Sikkim', 85, 39, 0, '26-06-2020', 15, 14, 0),
  @('Tamil Nadu', 70977, 39999, 911, '26-06-2020', 14132, 8683, 207),
  @('Telangana', 11364, 4688, 230, '26-06-2020', 4292, 1182, 27),
  @('Tripura', 1290, 1019, 1, '26-06-2020', 104, 339, 0),
  @('Uttarakhand', 2691, 1758, 36, '26-06-2020', 390, 308, 9),
  @('Uttar Pradesh', 20193, 13119, 611, '26-06-2020', 3599, 3124, 104),
  @('West Bengal', 15648, 10190, 606, '26-06-2020', 2117, 2325, 66)
)

$startRow = 3152
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $rowData = $newRows[$i]
  for ($j = 0; $j -lt $rowData.Count; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
  }
}

Write-Host "Added" $newRows.Count "rows; new UsedRange:" $ws.UsedRange.Address()
